# Auto-generated Excel COM-interop script
# Refreshes cached market-price / profit columns (H:N) on several rows
# across multiple sheets, per the scheduled market-data runner diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 380.65
$ws.Range("I28").Value = 448.75
$ws.Range("J28").Value = 278.5
$ws.Range("K28").Value = 448.75
$ws.Range("L28").Value = 278.5
$ws.Range("M28").Value = 36.25
$ws.Range("N28").Value = -1248.5
# Row 58
$ws.Range("H58").Value = 18001.826
$ws.Range("J58").Value = 21343.773
$ws.Range("L58").Value = 64031.319
$ws.Range("N58").Value = -64331.319
# Row 64
$ws.Range("H64").Value = 3028.5
$ws.Range("I64").Value = 3174.25
$ws.Range("J64").Value = 2931.3333
$ws.Range("K64").Value = 3174.25
$ws.Range("L64").Value = 2931.3333
$ws.Range("M64").Value = -2926.25
$ws.Range("N64").Value = -3427.3333
# Row 67
$ws.Range("H67").Value = 3028.5
$ws.Range("I67").Value = 3174.25
$ws.Range("J67").Value = 2931.3333
$ws.Range("K67").Value = 3174.25
$ws.Range("L67").Value = 2931.3333
$ws.Range("M67").Value = -2316.25
$ws.Range("N67").Value = -4647.3333
# Row 69
$ws.Range("H69").Value = 5115
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").Value = ""
# Row 72
$ws.Range("H72").Value = 5115
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").Value = ""
# Row 82
$ws.Range("H82").Value = 2783.6667
$ws.Range("I82").Value = 175.5
$ws.Range("K82").Value = 526.5
$ws.Range("M82").Value = -120.5
# Row 85
$ws.Range("H85").Value = 2783.6667
$ws.Range("I85").Value = 175.5
$ws.Range("K85").Value = 526.5
$ws.Range("M85").Value = 877.5
# Row 86
$ws.Range("H86").Value = 1383.3334
$ws.Range("I86").Value = 885.7143
$ws.Range("K86").Value = 885.7143
$ws.Range("M86").Value = 237.2857
# Row 89
$ws.Range("H89").Value = 1383.3334
$ws.Range("I89").Value = 885.7143
$ws.Range("K89").Value = 4428.5715
$ws.Range("M89").Value = 1187.4285
# Row 92
$ws.Range("H92").Value = 381.20834
$ws.Range("I92").Value = 382.45
$ws.Range("J92").Value = 375
$ws.Range("K92").Value = 382.45
$ws.Range("L92").Value = 375
$ws.Range("M92").Value = 865.55
$ws.Range("N92").Value = -2871
# Row 96
$ws.Range("H96").Value = 959.8
$ws.Range("I96").Value = 895
$ws.Range("J96").Value = 976
$ws.Range("K96").Value = 2685
$ws.Range("L96").Value = 2928
$ws.Range("M96").Value = -1312
$ws.Range("N96").Value = -5674
# Row 106
$ws.Range("H106").Value = 1003.5455
$ws.Range("I106").Value = 581.6667
$ws.Range("J106").Value = 1509.8
$ws.Range("K106").Value = 581.6667
$ws.Range("L106").Value = 1509.8
$ws.Range("M106").Value = 49.33330000000001
$ws.Range("N106").Value = -2771.8
# Row 115
$ws.Range("H115").Value = 1777
$ws.Range("I115").Value = 961.6667
$ws.Range("K115").Value = 2885.0001
$ws.Range("M115").Value = -1318.0001
# Row 129
$ws.Range("H129").Value = 3334493
$ws.Range("I129").Value = 35715468
$ws.Range("J129").Value = 1157.3677
$ws.Range("K129").Value = 107146404
$ws.Range("L129").Value = 3472.1031
$ws.Range("M129").Value = -107141404
$ws.Range("N129").Value = -13472.1031
# Row 135
$ws.Range("H135").Value = 1181.5161
$ws.Range("I135").Value = 870.2692
$ws.Range("J135").Value = 2800
$ws.Range("K135").Value = 7832.422799999999
$ws.Range("L135").Value = 25200
$ws.Range("M135").Value = -5297.422799999999
$ws.Range("N135").Value = -30270
# Row 137
$ws.Range("H137").Value = 2002533.5
$ws.Range("I137").Value = 3033067.2
$ws.Range("K137").Value = 9099201.600000001
$ws.Range("M137").Value = -9096651.600000001
# Row 141
$ws.Range("H141").Value = 218353.94
$ws.Range("I141").Value = 1088.1464
$ws.Range("J141").Value = 1490910.8
$ws.Range("K141").Value = 3264.4392
$ws.Range("L141").Value = 4472732.4
$ws.Range("M141").Value = 1915.5608
$ws.Range("N141").Value = -4483092.4

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3645.56
$ws.Range("I32").Value = 3057.976
$ws.Range("J32").Value = 6730.375
$ws.Range("K32").Value = 3057.976
$ws.Range("L32").Value = 6730.375
$ws.Range("M32").Value = -2770.976
$ws.Range("N32").Value = -7304.375
# Row 61
$ws.Range("H61").Value = 1702.4667
$ws.Range("I61").Value = 790.9429
$ws.Range("J61").Value = 4892.8
$ws.Range("K61").Value = 790.9429
$ws.Range("L61").Value = 4892.8
$ws.Range("M61").Value = -578.9429
$ws.Range("N61").Value = -5316.8
# Row 132
$ws.Range("H132").Value = 2196.7693
$ws.Range("I132").Value = 1685.9697
$ws.Range("K132").Value = 5057.909100000001
$ws.Range("M132").Value = -2527.909100000001
# Row 135
$ws.Range("H135").Value = 27809.166
$ws.Range("J135").Value = 27809.166
$ws.Range("L135").Value = 27809.166
$ws.Range("N135").Value = -37949.166
# Row 136
$ws.Range("H136").Value = 1702.4667
$ws.Range("I136").Value = 790.9429
$ws.Range("J136").Value = 4892.8
$ws.Range("K136").Value = 2372.8287
$ws.Range("L136").Value = 14678.4
$ws.Range("M136").Value = 177.1713
$ws.Range("N136").Value = -19778.4

$ws = $wb.Worksheets.Item("CRP")
# Row 13
$ws.Range("H13").Value = 67502.5
$ws.Range("J13").Value = 67502.5
$ws.Range("L13").Value = 67502.5
$ws.Range("N13").Value = -67780.5
# Row 107
$ws.Range("H107").Value = 2328.2666
$ws.Range("I107").Value = 701.375
$ws.Range("J107").Value = 4187.5713
$ws.Range("K107").Value = 701.375
$ws.Range("L107").Value = 4187.5713
$ws.Range("M107").Value = 1218.625
$ws.Range("N107").Value = -8027.5713
# Row 122
$ws.Range("H122").Value = 4218.6665
$ws.Range("I122").Value = 4632.3335
$ws.Range("J122").Value = 4149.722
$ws.Range("K122").Value = 13897.0005
$ws.Range("L122").Value = 12449.166
$ws.Range("M122").Value = -11447.0005
$ws.Range("N122").Value = -17349.166

$ws = $wb.Worksheets.Item("CUL")
# Row 63
$ws.Range("H63").Value = 3850
$ws.Range("J63").Value = 3971.4285
$ws.Range("L63").Value = 11914.2855
$ws.Range("N63").Value = -13412.2855
# Row 66
$ws.Range("H66").Value = 3850
$ws.Range("J66").Value = 3971.4285
$ws.Range("L66").Value = 35742.8565
$ws.Range("N66").Value = -43230.8565
# Row 87
$ws.Range("H87").Value = 5900.5557
$ws.Range("I87").Value = 2477.6924
$ws.Range("J87").Value = 14800
$ws.Range("K87").Value = 7433.0772
$ws.Range("L87").Value = 44400
$ws.Range("M87").Value = -6185.0772
$ws.Range("N87").Value = -46896
# Row 90
$ws.Range("H90").Value = 5900.5557
$ws.Range("I90").Value = 2477.6924
$ws.Range("J90").Value = 14800
$ws.Range("K90").Value = 22299.2316
$ws.Range("L90").Value = 133200
$ws.Range("M90").Value = -16059.2316
$ws.Range("N90").Value = -145680
# Row 120
$ws.Range("H120").Value = 15734.3
$ws.Range("I120").Value = 11077.5
$ws.Range("K120").Value = 33232.5
$ws.Range("M120").Value = -28394.5
# Row 134
$ws.Range("H134").Value = 2500.889
$ws.Range("I134").Value = 1036.6
$ws.Range("J134").Value = 4331.25
$ws.Range("K134").Value = 3109.8
$ws.Range("L134").Value = 12993.75
$ws.Range("M134").Value = 1960.2
$ws.Range("N134").Value = -23133.75
# Row 138
$ws.Range("H138").Value = 1738.9
$ws.Range("I138").Value = 903.1667
$ws.Range("J138").Value = 2992.5
$ws.Range("K138").Value = 2709.5001
$ws.Range("L138").Value = 8977.5
$ws.Range("M138").Value = 2430.4999
$ws.Range("N138").Value = -19257.5
# Row 140
$ws.Range("H140").Value = 16670647
$ws.Range("I140").Value = 41667390
$ws.Range("K140").Value = 125002170
$ws.Range("M140").Value = -124996990

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 18188.8
$ws.Range("J3").Value = 22236
$ws.Range("L3").Value = 22236
$ws.Range("N3").Value = -22468
# Row 97
$ws.Range("H97").Value = 1052.1578
$ws.Range("I97").Value = 427.14285
$ws.Range("K97").Value = 427.14285
$ws.Range("M97").Value = 68.85714999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 17
$ws.Range("H17").Value = 20000
$ws.Range("J17").Value = 20000
$ws.Range("L17").Value = 20000
$ws.Range("N17").Value = -20340
# Row 46
$ws.Range("H46").Value = 2035.7858
$ws.Range("I46").Value = 533.6667
$ws.Range("J46").Value = 2445.4546
$ws.Range("K46").Value = 533.6667
$ws.Range("L46").Value = 2445.4546
$ws.Range("M46").Value = -345.6667
$ws.Range("N46").Value = -2821.4546
# Row 68
$ws.Range("H68").Value = 2285.8215
$ws.Range("I68").Value = 1071.4286
$ws.Range("J68").Value = 5929
$ws.Range("K68").Value = 1071.4286
$ws.Range("L68").Value = 5929
$ws.Range("M68").Value = -322.4286
$ws.Range("N68").Value = -7427
# Row 71
$ws.Range("H71").Value = 2285.8215
$ws.Range("I71").Value = 1071.4286
$ws.Range("J71").Value = 5929
$ws.Range("K71").Value = 5357.143
$ws.Range("L71").Value = 29645
$ws.Range("M71").Value = -1613.143
$ws.Range("N71").Value = -37133
# Row 112
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = ""
# Row 136
$ws.Range("H136").Value = 4004059.8
$ws.Range("I136").Value = 5884796.5
$ws.Range("K136").Value = 17654389.5
$ws.Range("M136").Value = -17651839.5

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 252774.33
$ws.Range("I132").Value = 347793.4
$ws.Range("J132").Value = 40808.69
$ws.Range("K132").Value = 1043380.2
$ws.Range("L132").Value = 122426.07
$ws.Range("M132").Value = -1040850.2
$ws.Range("N132").Value = -127486.07
